$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from an existing header
# cell onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF), rows 2-17
$iValues = @(7, 9, 9, 9, 9, 7, 7, 7, 8, 8, 8, 5, 6, 1, 3, 3)
$jValues = @(8, 9, 9, 9, 9, 7, 8, 8, 9, 9, 8, 7, 8, 2, 4, 3)

for ($r = 0; $r -lt 16; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
